# Rebuild the Fatty-Acid whitelist table: rename header sn1/sn2/sn3 -> fa1/fa2/fa3,
# and splice in newly-whitelisted short/odd-chain fatty acids (FA4:0..FA11:0, FA25:0)
# in numeric order, growing the table from A1:G55 to A1:G62. Fixes the TG [M+Na]+ /
# [M+H]+ control bug.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -- introduce the brand-new labels first (matches author's original edit order) --
$ws.Range("B1").Value = "fa1"
$ws.Range("C1").Value = "fa2"
$ws.Range("D1").Value = "fa3"
$ws.Range("A9").Value = "FA11:0"
$ws.Range("A2").Value = "FA4:0"
$ws.Range("A3").Value = "FA5:0"
$ws.Range("A4").Value = "FA6:0"
$ws.Range("A5").Value = "FA7:0"
$ws.Range("A7").Value = "FA9:0"
$ws.Range("A53").Value = "FA25:0"

# -- now lay down the full A1:G62 table (re-asserts the cells set above too) --
$data = @(
    ,@("FattyAcid", "fa1", "fa2", "fa3", "PL", "TG", "DG")
    ,@("FA4:0", "T", "T", "T", $null, "T", "T")
    ,@("FA5:0", "T", "T", "T", $null, "T", "T")
    ,@("FA6:0", "T", "T", "T", $null, "T", "T")
    ,@("FA7:0", "T", "T", "T", $null, "T", "T")
    ,@("FA8:0", "T", "T", "T", $null, "T", "T")
    ,@("FA9:0", "T", "T", "T", $null, "T", "T")
    ,@("FA10:0", "T", "T", "T", $null, "T", "T")
    ,@("FA11:0", "T", "T", "T", $null, "T", "T")
    ,@("FA12:0", "T", "T", "T", $null, "T", "T")
    ,@("FA12:1", "T", "T", "T", $null, "T", "T")
    ,@("FA13:0", "T", "T", "T", $null, "T", "T")
    ,@("FA13:1", "T", "T", "T", $null, "T", "T")
    ,@("FA13:2", "T", "T", "T", $null, "T", "T")
    ,@("FA14:0", "T", "T", "T", $null, "T", "T")
    ,@("FA14:1", "T", "T", "T", $null, "T", "T")
    ,@("FA14:2", "T", "T", "T", $null, "T", "T")
    ,@("FA15:0", "T", "T", "T", $null, "T", "T")
    ,@("FA15:1", "T", "T", "T", $null, "T", "T")
    ,@("FA15:2", "T", "T", "T", $null, "T", "T")
    ,@("FA16:0", "T", "T", "T", "T", "T", "T")
    ,@("FA16:1", "T", "T", "T", $null, "T", "T")
    ,@("FA16:2", "T", "T", "T", $null, "T", "T")
    ,@("FA17:0", "T", "T", "T", $null, "T", "T")
    ,@("FA17:1", "T", "T", "T", $null, "T", "T")
    ,@("FA17:2", "T", "T", "T", $null, "T", "T")
    ,@("FA18:0", "T", "T", "T", "T", "T", "T")
    ,@("FA18:1", "T", "T", "T", "T", "T", "T")
    ,@("FA18:2", "T", "T", "T", "T", "T", "T")
    ,@("FA18:3", "T", "T", "T", "T", "T", "T")
    ,@("FA18:4", "T", "T", "T", $null, "T", "T")
    ,@("FA19:0", "T", "T", "T", $null, "T", "T")
    ,@("FA19:1", "T", "T", "T", $null, "T", "T")
    ,@("FA19:2", "T", "T", "T", $null, "T", "T")
    ,@("FA20:0", "T", "T", "T", $null, "T", "T")
    ,@("FA20:1", "T", "T", "T", $null, "T", "T")
    ,@("FA20:2", "T", "T", "T", $null, "T", "T")
    ,@("FA20:3", "T", "T", "T", "T", "T", "T")
    ,@("FA20:4", "T", "T", "T", "T", "T", "T")
    ,@("FA20:5", "T", "T", "T", "T", "T", "T")
    ,@("FA22:0", "T", "T", "T", $null, "T", "T")
    ,@("FA22:1", "T", "T", "T", $null, "T", "T")
    ,@("FA22:2", "T", "T", "T", $null, "T", "T")
    ,@("FA22:3", "T", "T", "T", $null, "T", "T")
    ,@("FA22:4", "T", "T", "T", "T", "T", "T")
    ,@("FA22:5", "T", "T", "T", "T", "T", "T")
    ,@("FA22:6", "T", "T", "T", "T", "T", $null)
    ,@("FA23:0", "T", "T", "T", $null, "T", $null)
    ,@("FA24:0", "T", "T", "T", $null, "T", $null)
    ,@("FA24:1", "T", "T", "T", $null, "T", $null)
    ,@("FA24:2", "T", "T", "T", $null, "T", $null)
    ,@("FA24:3", "T", "T", "T", $null, "T", $null)
    ,@("FA25:0", "T", "T", "T", $null, "T", $null)
    ,@("FA26:0", "T", "T", "T", $null, "T", $null)
    ,@("FA26:1", "T", "T", "T", $null, "T", $null)
    ,@("FA26:2", "T", "T", "T", $null, "T", $null)
    ,@("O-16:0", $null, $null, $null, "T", $null, $null)
    ,@("O-18:0", "T", $null, $null, "T", $null, $null)
    ,@("O-20:0", "T", $null, $null, "T", $null, $null)
    ,@("P-16:0", "T", $null, $null, "T", $null, $null)
    ,@("P-18:0", "T", $null, $null, "T", $null, $null)
    ,@("P-20:0", "T", $null, $null, "T", $null, $null)
)

$r = 1
foreach ($row in $data) {
    $c = 1
    foreach ($val in $row) {
        if ($val -eq $null) {
            $ws.Cells.Item($r, $c).ClearContents()
        } else {
            $ws.Cells.Item($r, $c).Value = $val
        }
        $c = $c + 1
    }
    $r = $r + 1
}

# -- match the author's final scroll position / selection in the saved view state --
$win = $excel.ActiveWindow
try { $win.TopLeftCell = $ws.Range("A16") } catch { }
$ws.Range("F53").Select()
